$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DataSet")
$ws.Rows.Item(40).Delete()
